# EPBDS: added validation phase for mapping subproject.
# Sheet1 ("Data Mapping mappings1"): the D18:D22 column ("aString" field name)
# used to repeat the same text in every row. Turn it into a single merged
# cell D18:D22 (matching the existing C18:C22 / F18:F22 / H18:H22 / I18:I22
# merges), keeping the "aString" label only in the top cell and clearing
# the duplicated text from D19:D22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D19:D22 should end up blank, but keep the centered look already used by
# the neighboring C19:C22 cells (style s="16": horizontal/vertical center).
$ws.Range("D19:D22").ClearContents()

# Merge D18:D22 into a single cell, just like the other columns in this block.
$ws.Range("D18:D22").Merge()

# D18 should look like the other merged header cells on row 18 (centered,
# with the thin top border carried by style s="15"); F18 already has that
# exact formatting, so copy formats only (xlPasteFormats = -4122).
$ws.Range("F18").Copy()
$ws.Range("D18").PasteSpecial(-4122)

# D19:D22 should look like C19:C22 (style s="16": centered, no border).
$ws.Range("C19").Copy()
$ws.Range("D19:D22").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Move the active selection to the newly merged cell, like in the edit.
$ws.Range("D18:D22").Select() | Out-Null

Write-Output "Merged D18:D22 and refreshed formatting"
